# Ajout sélecteur meilleur candidat.
# Recompute the "Ordre" (play-order / best-candidate rank) column (H) on the
# "Data" sheet of the ListeChansons table, reflecting the new best-candidate
# selection order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 1
$ws.Range("H3").Value  = 37
$ws.Range("H4").Value  = 16
$ws.Range("H5").Value  = 32
$ws.Range("H6").Value  = 39
$ws.Range("H7").Value  = 26
$ws.Range("H8").Value  = 31
$ws.Range("H9").Value  = 35
$ws.Range("H10").Value = 38
$ws.Range("H11").Value = 36
$ws.Range("H12").Value = 24
$ws.Range("H13").Value = 40
$ws.Range("H14").Value = 21
$ws.Range("H15").Value = 12
$ws.Range("H16").Value = 4
$ws.Range("H17").Value = 13
$ws.Range("H18").Value = 14
$ws.Range("H19").Value = 28
$ws.Range("H20").Value = 22
$ws.Range("H21").Value = 20
$ws.Range("H22").Value = 23
$ws.Range("H23").Value = 2
$ws.Range("H24").Value = 17
$ws.Range("H25").Value = 9
$ws.Range("H26").Value = 6
$ws.Range("H27").Value = 3
$ws.Range("H28").Value = 7
$ws.Range("H29").Value = 25
$ws.Range("H30").Value = 30
$ws.Range("H31").Value = 11
$ws.Range("H32").Value = 27
$ws.Range("H33").Value = 29
$ws.Range("H34").Value = 8
$ws.Range("H35").Value = 15
$ws.Range("H36").Value = 33
$ws.Range("H37").Value = 10
$ws.Range("H38").Value = 18
$ws.Range("H39").Value = 34
$ws.Range("H40").Value = 19
$ws.Range("H41").Value = 5
